$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 updates ---
$ws.Range("P6").Value = 61
$ws.Range("Q6").Value = "Charter Privatflugzeug"
$ws.Range("V6").Value = "Charterflug Privat"
$ws.Range("AD6").Value = 3
$ws.Range("AG6").Value = "Privatflugzeug"

$ws.Range("AI6").Value = "2341-M"
$ws.Range("AJ6").Value = 5
$ws.Range("AL6").Value = "EXW-FS"
$ws.Range("AM6").Value = 12
$ws.Range("AO6").Value = 6203
$ws.Range("AP6").Value = 1
$ws.Range("AQ6").Value = "Landung"

# Remove the now-unused trailing cells (shift structure left, not rows)
$ws.Range("AR6:AT6").Clear()

# --- Row 11 updates ---
$ws.Range("P11").Value = 61
$ws.Range("Q11").Value = "Charter Privatflugzeug"
$ws.Range("V11").Value = "Charterflug Privat"
$ws.Range("AD11").Value = 2
$ws.Range("AG11").Value = "Privatflugzeug"

$ws.Range("AI11").Value = "2341-M"
$ws.Range("AJ11").Value = 3
$ws.Range("AL11").Value = "EXW-FS-U"
$ws.Range("AM11").Value = 25

# Remove the now-unused trailing cells
$ws.Range("AO11:AQ11").Clear()

# --- Sheet view state changes ---
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("A2").Select()
$ws.Range("AD7").Select()
